# Applies the edits described by the commit "fix erro OS desc ass de negocios"
$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

ReplaceText "NOME: RENAN NUNES" "NOME: AAAAA"
ReplaceText "FUNÇÃO: ANALISTA DE CUSTO" "FUNÇÃO: ANALISTA DE NEGOCIOS"
ReplaceText "MECÂNICO" "ACIDENTE"
ReplaceText "RENAN NUNES" "AAAAA"
ReplaceText "ANALISTA DE CUSTO" "ANALISTA DE NEGOCIOS"
ReplaceText "MANOEL JEFETE DA SILVA TENONIO" "BRUNA PETRONI CEZARIO"
ReplaceText "MTE/RN: 1805" "CREA-RN: 2122993685"
ReplaceText "Parnamirim/RN, 30 de abril de 2025." "Parnamirim/RN, 30 de Junho de 2025."
